# Automatische test-sync: 2025-08-05 18:28:50
# Appends a new test-mail log entry (row 32) to the "Logs" sheet and
# updates the "Dashboard" summary sheet to reflect the new category count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New row 32 on the Logs sheet -----------------------------------------
$ws.Range("A32").Value = "Mijn retour is nog steeds niet verwerkt."
$ws.Range("B32").Value = "mailmind.test@zohomail.eu"
$ws.Range("C32").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$ws.Range("D32").Value = "Retour / Terugbetaling"
$ws.Range("E32").Value = "Beste klant,`r`nBedankt voor je bericht. Ik begrijp dat je retour nog niet is verwerkt en ik wil je graag helpen om dit op te lossen.`r`nOm dit verder te onderzoeken, heb ik wat meer informatie nodig. Zou je zo vriendelijk willen zijn om je ordernummer en/of trackingnummer met me te delen? Hiermee kan ik de status van je retour nakijken en je zo goed mogelijk van dienst zijn.`r`nIk kijk uit naar je reactie.`r`nMet vriendelijke groet,`r`n[Naam van de e-mailassistent]  `r`nJamie  `r`nNederlandse e-mailassistent  `r`n[Bedrijfsnaam]"
$ws.Range("F32").Value = "2025-08-05 18:28:24"
$ws.Range("G32").Value = "Ja"
$ws.Range("H32").Value = "Nee"
$ws.Range("I32").Value = "Ja"
$ws.Range("J32").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row -----
$ws.Range("D2:D31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D32"))
$ws.Range("G2:G31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G32"))
$ws.Range("H2:H31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H32"))
$ws.Range("I2:I31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I32"))
$ws.Range("J2:J31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J32"))

# --- Update the Dashboard summary table ------------------------------------
# "Retour / Terugbetaling" now has 2 entries (was 1) and moves ahead of
# "Opvolging / Status" (still 2), which drops to the row below.
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("A5").Value = "Retour / Terugbetaling"
$wsDash.Range("B5").Value = 2
$wsDash.Range("A6").Value = "Opvolging / Status"
$wsDash.Range("B6").Value = 2
